# Add new session 28-07-23
# Appends two new data rows (8 and 9) to the bottom of the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 246077
$ws.Range("B8").Value = "random"
$ws.Range("C8").Value = "2023-07-21T15:19:19.000000Z"
$ws.Range("D8").Value = 9.720000000000001
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 30
$ws.Range("J8").Value = 8
$ws.Range("K8").Value = 14
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 4
$ws.Range("N8").Value = 11
$ws.Range("O8").Value = 3
$ws.Range("P8").Value = 11
$ws.Range("Q8").Value = 21
$ws.Range("R8").Value = 3
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 15
$ws.Range("U8").Value = 14
$ws.Range("V8").Value = 2
$ws.Range("W8").Value = 6
$ws.Range("X8").Value = 7
$ws.Range("Y8").Value = 9

# Row 9
$ws.Range("A9").Value = 248429
$ws.Range("B9").Value = "random"
$ws.Range("C9").Value = "2023-07-28T18:38:43.000000Z"
$ws.Range("D9").Value = 9.25
$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 12
$ws.Range("I9").Value = 19
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 17
$ws.Range("L9").Value = 13
$ws.Range("M9").Value = 17
$ws.Range("N9").Value = 16
$ws.Range("O9").Value = 15
$ws.Range("P9").Value = 4
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 6
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 5
$ws.Range("U9").Value = 23
$ws.Range("V9").Value = 19
$ws.Range("W9").Value = 22
$ws.Range("X9").Value = 9
$ws.Range("Y9").Value = 5
